$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 955.8570999999999
$ws.Range("I19").Value = 900.5
$ws.Range("J19").Value = 978
$ws.Range("K19").Value = 900.5
$ws.Range("L19").Value = 978
$ws.Range("M19").Value = -725.5
$ws.Range("N19").Value = -1328
$ws.Range("H98").Value = 267.66666
$ws.Range("I98").Value = 215.35715
$ws.Range("K98").Value = 215.35715
$ws.Range("M98").Value = 1282.64285
$ws.Range("H112").Value = 3969281
$ws.Range("J112").Value = 4274564.5
$ws.Range("L112").Value = 12823693.5
$ws.Range("N112").Value = -12825909.5
$ws.Range("H122").Value = 267.66666
$ws.Range("I122").Value = 215.35715
$ws.Range("K122").Value = 646.0714499999999
$ws.Range("M122").Value = 1803.92855
$ws.Range("H138").Value = 1663.963
$ws.Range("I138").Value = 610.5625
$ws.Range("K138").Value = 1831.6875
$ws.Range("M138").Value = 3308.3125
$ws.Range("H141").Value = 3050
$ws.Range("I141").Value = 1700
$ws.Range("J141").Value = 4400
$ws.Range("K141").Value = 5100
$ws.Range("L141").Value = 13200
$ws.Range("M141").Value = 80
$ws.Range("N141").Value = -23560

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1762.7
$ws.Range("J2").Value = 3237.6667
$ws.Range("L2").Value = 3237.6667
$ws.Range("N2").Value = -3463.6667
$ws.Range("H24").Value = 26000
$ws.Range("J24").Value = 26000
$ws.Range("L24").Value = 26000
$ws.Range("N24").Value = -26748
$ws.Range("H32").Value = 43636.96
$ws.Range("I32").Value = 45731.543
$ws.Range("J32").Value = 18502
$ws.Range("K32").Value = 45731.543
$ws.Range("L32").Value = 18502
$ws.Range("M32").Value = -45444.543
$ws.Range("N32").Value = -19076
$ws.Range("H61").Value = 1622.6222
$ws.Range("I61").Value = 1004.5294
$ws.Range("K61").Value = 1004.5294
$ws.Range("M61").Value = -792.5294
$ws.Range("H74").Value = 2611.963
$ws.Range("I74").Value = 2724.158
$ws.Range("J74").Value = 2345.5
$ws.Range("K74").Value = 2724.158
$ws.Range("L74").Value = 2345.5
$ws.Range("M74").Value = -1850.158
$ws.Range("N74").Value = -4093.5
$ws.Range("H77").Value = 2611.963
$ws.Range("I77").Value = 2724.158
$ws.Range("J77").Value = 2345.5
$ws.Range("K77").Value = 13620.79
$ws.Range("L77").Value = 11727.5
$ws.Range("M77").Value = -9252.789999999999
$ws.Range("N77").Value = -20463.5
$ws.Range("H88").Value = 201190
$ws.Range("I88").Value = 1748
$ws.Range("J88").Value = 334151.34
$ws.Range("K88").Value = 1748
$ws.Range("L88").Value = 334151.34
$ws.Range("M88").Value = -1342
$ws.Range("N88").Value = -334963.34
$ws.Range("H91").Value = 201190
$ws.Range("I91").Value = 1748
$ws.Range("J91").Value = 334151.34
$ws.Range("K91").Value = 1748
$ws.Range("L91").Value = 334151.34
$ws.Range("M91").Value = -344
$ws.Range("N91").Value = -336959.34
$ws.Range("H100").Value = 26000
$ws.Range("J100").Value = 26000
$ws.Range("L100").Value = 26000
$ws.Range("N100").Value = -28164
$ws.Range("H102").Value = 1249.875
$ws.Range("I102").Value = 1115.2307
$ws.Range("J102").Value = 1833.3334
$ws.Range("K102").Value = 1115.2307
$ws.Range("L102").Value = 1833.3334
$ws.Range("M102").Value = 506.7692999999999
$ws.Range("N102").Value = -5077.3334
$ws.Range("H110").Value = 389.3125
$ws.Range("I110").Value = 349.93332
$ws.Range("J110").Value = 980
$ws.Range("K110").Value = 349.93332
$ws.Range("L110").Value = 980
$ws.Range("M110").Value = 1695.06668
$ws.Range("N110").Value = -5070
$ws.Range("H116").Value = 1762.7
$ws.Range("J116").Value = 3237.6667
$ws.Range("L116").Value = 3237.6667
$ws.Range("N116").Value = -7825.6667
$ws.Range("H122").Value = 2226.8462
$ws.Range("I122").Value = 2268.1365
$ws.Range("J122").Value = 1999.75
$ws.Range("K122").Value = 6804.4095
$ws.Range("L122").Value = 5999.25
$ws.Range("M122").Value = -4354.4095
$ws.Range("N122").Value = -10899.25
$ws.Range("H136").Value = 1622.6222
$ws.Range("I136").Value = 1004.5294
$ws.Range("K136").Value = 3013.5882
$ws.Range("M136").Value = -463.5882000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1762.7
$ws.Range("J3").Value = 3237.6667
$ws.Range("L3").Value = 3237.6667
$ws.Range("N3").Value = -3465.6667
$ws.Range("H86").Value = 1648.0714
$ws.Range("I86").Value = 1386.25
$ws.Range("J86").Value = 3219
$ws.Range("K86").Value = 1386.25
$ws.Range("L86").Value = 3219
$ws.Range("M86").Value = -263.25
$ws.Range("N86").Value = -5465
$ws.Range("H89").Value = 1648.0714
$ws.Range("I89").Value = 1386.25
$ws.Range("J89").Value = 3219
$ws.Range("K89").Value = 6931.25
$ws.Range("L89").Value = 16095
$ws.Range("M89").Value = -1315.25
$ws.Range("N89").Value = -27327
$ws.Range("H99").Value = 1959.0555
$ws.Range("I99").Value = 1661.1818
$ws.Range("J99").Value = 2427.1428
$ws.Range("K99").Value = 1661.1818
$ws.Range("L99").Value = 2427.1428
$ws.Range("M99").Value = -163.1818000000001
$ws.Range("N99").Value = -5423.1428
$ws.Range("H107").Value = 1189.4706
$ws.Range("I107").Value = 392.5
$ws.Range("J107").Value = 3102.2
$ws.Range("K107").Value = 392.5
$ws.Range("L107").Value = 3102.2
$ws.Range("M107").Value = 1527.5
$ws.Range("N107").Value = -6942.2
$ws.Range("H134").Value = 38866.965
$ws.Range("I134").Value = 47033.086
$ws.Range("J134").Value = 1302.8
$ws.Range("K134").Value = 141099.258
$ws.Range("L134").Value = 3908.4
$ws.Range("M134").Value = -138564.258
$ws.Range("N134").Value = -8978.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 16670200
$ws.Range("I99").Value = 3971035
$ws.Range("K99").Value = 3971035
$ws.Range("M99").Value = -3969537
$ws.Range("H107").Value = 1131.3636
$ws.Range("I107").Value = 710
$ws.Range("J107").Value = 1740
$ws.Range("K107").Value = 710
$ws.Range("L107").Value = 1740
$ws.Range("M107").Value = 1210
$ws.Range("N107").Value = -5580
$ws.Range("H126").Value = 16670200
$ws.Range("I126").Value = 3971035
$ws.Range("K126").Value = 11913105
$ws.Range("M126").Value = -11910635

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 3253.3547
$ws.Range("I2").Value = 4372.696
$ws.Range("K2").Value = 26236.176
$ws.Range("M2").Value = -26123.176
$ws.Range("H13").Value = 188
$ws.Range("I13").Value = 132
$ws.Range("K13").Value = 396
$ws.Range("M13").Value = -228
$ws.Range("H14").Value = 282
$ws.Range("I14").Value = 282
$ws.Range("K14").Value = 846
$ws.Range("M14").Value = -673
$ws.Range("H36").Value = 2899.8333
$ws.Range("I36").Value = 2931.3333
$ws.Range("J36").Value = 2868.3333
$ws.Range("K36").Value = 8793.999899999999
$ws.Range("L36").Value = 8604.999899999999
$ws.Range("M36").Value = -8624.999899999999
$ws.Range("N36").Value = -8942.999899999999
$ws.Range("H95").Value = 5027
$ws.Range("J95").Value = 5027
$ws.Range("L95").Value = 15081
$ws.Range("N95").Value = -19199
$ws.Range("H131").Value = 761.0700000000001
$ws.Range("J131").Value = 772.57446
$ws.Range("L131").Value = 2317.72338
$ws.Range("N131").Value = -12397.72338

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 16877
$ws.Range("I44").Value = 18000
$ws.Range("J44").Value = 16315.5
$ws.Range("K44").Value = 18000
$ws.Range("L44").Value = 16315.5
$ws.Range("M44").Value = -17404
$ws.Range("N44").Value = -17507.5
$ws.Range("H102").Value = 2230.375
$ws.Range("I102").Value = 2230.375
$ws.Range("K102").Value = 2230.375
$ws.Range("M102").Value = -608.375
$ws.Range("H126").Value = 4251.032
$ws.Range("I126").Value = 3152.2104
$ws.Range("K126").Value = 9456.6312
$ws.Range("M126").Value = -6986.6312

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1776.5454
$ws.Range("I61").Value = 1407
$ws.Range("J61").Value = 2568.4285
$ws.Range("K61").Value = 1407
$ws.Range("L61").Value = 2568.4285
$ws.Range("M61").Value = -1205
$ws.Range("N61").Value = -2972.4285
$ws.Range("H94").Value = 35000
$ws.Range("J94").Value = 35000
$ws.Range("L94").Value = 35000
$ws.Range("N94").Value = -36352
$ws.Range("H100").Value = 1966.5555
$ws.Range("I100").Value = 1754.8182
$ws.Range("K100").Value = 1754.8182
$ws.Range("M100").Value = -1213.8182
$ws.Range("H101").Value = 13021.5
$ws.Range("J101").Value = 13021.5
$ws.Range("L101").Value = 13021.5
$ws.Range("N101").Value = -19511.5
$ws.Range("H113").Value = 1776.5454
$ws.Range("I113").Value = 1407
$ws.Range("J113").Value = 2568.4285
$ws.Range("K113").Value = 1407
$ws.Range("L113").Value = 2568.4285
$ws.Range("M113").Value = 763
$ws.Range("N113").Value = -6908.4285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1567908.2
$ws.Range("I107").Value = 514.5
$ws.Range("K107").Value = 1543.5
$ws.Range("M107").Value = 376.5
$ws.Range("H113").Value = 1229232
$ws.Range("I113").Value = 848.875
$ws.Range("J113").Value = 4504920
$ws.Range("K113").Value = 2546.625
$ws.Range("L113").Value = 13514760
$ws.Range("M113").Value = -376.625
$ws.Range("N113").Value = -13519100
$ws.Range("H122").Value = 1833.6522
$ws.Range("I122").Value = 1866.25
$ws.Range("J122").Value = 1616.3334
$ws.Range("K122").Value = 5598.75
$ws.Range("L122").Value = 4849.0002
$ws.Range("M122").Value = -3148.75
$ws.Range("N122").Value = -9749.0002
$ws.Range("H126").Value = 1504.3889
$ws.Range("I126").Value = 1013.1667
$ws.Range("J126").Value = 1750
$ws.Range("K126").Value = 3039.5001
$ws.Range("L126").Value = 5250
$ws.Range("M126").Value = -569.5001000000002
$ws.Range("N126").Value = -10190
